# updated statbar xpaths & diagnosis testcases
#
# Adds three new worksheets (CypherOutput_Message, StatOutput,
# StatOutput_Message): CypherOutput_Message mirrors the existing
# "Message" sheet, StatOutput carries the new statbar counts, and
# StatOutput_Message repeats the message block twice, the second
# occurrence using the updated (statbar) Cypher query.

$wb = $excel.ActiveWorkbook

$messageLines = @(
    'Neo4j_URL:',
    'bolt://ncias-q2251-c.nci.nih.gov:7687',
    'User_name:',
    'neo4j',
    'PWD:',
    'icdcDBneo4j0',
    'Cypher:',
    'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Lip and oropharyngeal neoplasms malignant :: Melanoma-lingual''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`',
    'Output:',
    'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC02_Canine_Filter_Diagnosis-LipLingual_Neo4jData.xlsx'
)

$statCypher = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Lip and oropharyngeal neoplasms malignant :: Melanoma-lingual'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# ---------------------------------------------------------------------
# 1) CypherOutput_Message - exact copy of the "Message" sheet content.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$cypherOutputMessage.Name = "CypherOutput_Message"

for ($i = 0; $i -lt $messageLines.Length; $i++) {
    $cypherOutputMessage.Cells.Item($i + 1, 1).Value = $messageLines[$i]
}

# ---------------------------------------------------------------------
# 2) StatOutput - summary counts table (headers + one data row).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutput = $wb.Worksheets.Add($null, $lastSheet)
$statOutput.Name = "StatOutput"

$statOutput.Cells.Item(1, 1).Value = "number_of_files"
$statOutput.Cells.Item(1, 2).Value = "number_of_sample"
$statOutput.Cells.Item(1, 3).Value = "number_of_cases"
$statOutput.Cells.Item(1, 4).Value = "number_of_study"

# Counts are stored as text (not numbers) in the source report, so
# force the Text number format before writing the digit strings --
# otherwise Excel auto-converts them to numeric cells.
$statOutput.Range("A2:D2").NumberFormat = "@"
$statOutput.Cells.Item(2, 1).Value = "19"
$statOutput.Cells.Item(2, 2).Value = "2"
$statOutput.Cells.Item(2, 3).Value = "1"
$statOutput.Cells.Item(2, 4).Value = "1"

# ---------------------------------------------------------------------
# 3) StatOutput_Message - the "Message" block repeated twice; the
#    second block's Cypher line is the updated statbar query.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$statOutputMessage.Name = "StatOutput_Message"

for ($i = 0; $i -lt $messageLines.Length; $i++) {
    $statOutputMessage.Cells.Item($i + 1, 1).Value = $messageLines[$i]
}

for ($i = 0; $i -lt $messageLines.Length; $i++) {
    $value = $messageLines[$i]
    if ($i -eq 7) {
        $value = $statCypher
    }
    $statOutputMessage.Cells.Item(10 + $i + 1, 1).Value = $value
}

# Restore the original active tab (adding sheets makes the newest one
# active, but the source workbook keeps "CypherOutput" selected).
$wb.Worksheets.Item("CypherOutput").Activate()

Write-Output "done"
